$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4: columns B, D, E, F decrease by 1 (6->5, 6->5, 6->5, 6->4)
foreach ($r in 2..4) {
    $ws.Range("B$r").Value = 5
    $ws.Range("D$r").Value = 5
    $ws.Range("E$r").Value = 5
    $ws.Range("F$r").Value = 4
}

# Rows 5-7: columns B, D, F decrease by 1 (6->5)
foreach ($r in 5..7) {
    $ws.Range("B$r").Value = 5
    $ws.Range("D$r").Value = 5
    $ws.Range("F$r").Value = 5
}

# Rows 8-20: columns B, D, F decrease by 1 (7->6)
foreach ($r in 8..20) {
    $ws.Range("B$r").Value = 6
    $ws.Range("D$r").Value = 6
    $ws.Range("F$r").Value = 6
}

# Update the sheet view: remove topLeftCell="A2" (scroll back to top-left)
# and change the selected cell from G21 to E4.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E4").Select()
